# Generate Report for Handoff
#
# This refreshes the localization-status report: the two files that were
# previously "Handed back: in sync with en-US" (28872837-...md and
# 8632ade2-...md) are now reported as "Ready for handoff" again, with
# refreshed handoff timestamps and (for zh-cn / de-de) a new error detail
# explaining that the handback file on record is stale.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("E2").Value = $readyForHandoff
$ovw.Range("F2").Value = $readyForHandoff
$ovw.Range("G2").Value = "2016-08-21 14:35:48"

$ovw.Range("E3").Value = $readyForHandoff
$ovw.Range("F3").Value = $readyForHandoff
$ovw.Range("G3").Value = "2016-08-21 14:35:48"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $readyForHandoff
$zhcn.Range("H2").Value = "2016-08-21 14:35:43"
$zhcn.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea775b8257240a33c46ec7249699a49718eae5c5/e2e/28872837-25ca-4f06-ab42-aaeffbb12d29.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/663dfffd3432ef76801bafc04f8d81a791b91adb/e2e/28872837-25ca-4f06-ab42-aaeffbb12d29.md."

$zhcn.Range("C3").Value = $readyForHandoff
$zhcn.Range("H3").Value = "2016-08-21 14:35:43"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea775b8257240a33c46ec7249699a49718eae5c5/e2e/8632ade2-8357-457b-a3ba-10e439bb9edf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/663dfffd3432ef76801bafc04f8d81a791b91adb/e2e/8632ade2-8357-457b-a3ba-10e439bb9edf.md."

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $readyForHandoff
$dede.Range("H2").Value = "2016-08-21 14:35:48"
$dede.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea775b8257240a33c46ec7249699a49718eae5c5/e2e/28872837-25ca-4f06-ab42-aaeffbb12d29.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/663dfffd3432ef76801bafc04f8d81a791b91adb/e2e/28872837-25ca-4f06-ab42-aaeffbb12d29.md."

$dede.Range("C3").Value = $readyForHandoff
$dede.Range("H3").Value = "2016-08-21 14:35:48"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea775b8257240a33c46ec7249699a49718eae5c5/e2e/8632ade2-8357-457b-a3ba-10e439bb9edf.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/663dfffd3432ef76801bafc04f8d81a791b91adb/e2e/8632ade2-8357-457b-a3ba-10e439bb9edf.md."

# ---------------------------------------------------------------------
# Column width adjustments (auto-fit side effect of the shorter/longer
# status & error text now in these columns). ColumnWidth is specified in
# characters (~1/6 increments); the values below are the closest
# achievable match to the target OOXML column widths of ~17.216 and 40.
# ---------------------------------------------------------------------
$narrowStatusColumnWidth = 16.333333333333336
$wideErrorColumnWidth = 39.166666666666664

$ovw.Range("E:F").ColumnWidth = $narrowStatusColumnWidth

$zhcn.Range("C:C").ColumnWidth = $narrowStatusColumnWidth
$zhcn.Range("P:P").ColumnWidth = $wideErrorColumnWidth

$dede.Range("C:C").ColumnWidth = $narrowStatusColumnWidth
$dede.Range("P:P").ColumnWidth = $wideErrorColumnWidth
